$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.61"
$ws.Range("E2").Value = "'0.16%"
$ws.Range("D3").Value = "'44.18"
$ws.Range("E3").Value = "'-0.12%"
$ws.Range("D4").Value = "'5.508"
$ws.Range("E4").Value = "'-1.01%"
$ws.Range("D5").Value = "'0.08080"
$ws.Range("E5").Value = "'0.24%"
$ws.Range("D6").Value = "'2.024"
$ws.Range("E6").Value = "'6.00%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9516"
$ws.Range("E7").Value = "'0.11%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1132"
$ws.Range("E8").Value = "'-6.67%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1873"
$ws.Range("E9").Value = "'1.61%"
$ws.Range("B10").Value = "MCDex"
$ws.Range("C10").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D10").Value = "'10.08"
$ws.Range("E10").Value = "'0.71%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09974"
$ws.Range("E11").Value = "'2.93%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04838"
$ws.Range("E12").Value = "'10.86%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1057"
$ws.Range("E13").Value = "'-0.63%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001254"
$ws.Range("E14").Value = "'-2.52%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04081"
$ws.Range("E15").Value = "'-3.18%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006003"
$ws.Range("E16").Value = "'1.00%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.372"
$ws.Range("E17").Value = "'-0.66%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.410"
$ws.Range("E18").Value = "'2.99%"
$ws.Range("D19").Value = "'2.621"
$ws.Range("E19").Value = "'1.78%"
$ws.Range("D20").Value = "'0.3295"
$ws.Range("E20").Value = "'-4.93%"
$ws.Range("D21").Value = "'0.1399"
$ws.Range("E21").Value = "'-1.72%"
$ws.Range("D22").Value = "'0.2572"
$ws.Range("E22").Value = "'2.73%"
$ws.Range("D23").Value = "'0.001307"
$ws.Range("E23").Value = "'5.06%"
$ws.Range("D24").Value = "'0.004364"
$ws.Range("E24").Value = "'1.38%"
$ws.Range("E25").Value = "'4.94%"
$ws.Range("D26").Value = "'0.0003739"
$ws.Range("D38").Value = "'0.02589"
$ws.Range("E38").Value = "'-3.35%"
$ws.Range("D39").Value = "'0.05632"
$ws.Range("E39").Value = "'2.02%"
$ws.Range("D40").Value = "'0.007622"
$ws.Range("E40").Value = "'0.21%"
$ws.Range("D41").Value = "'0.1399"
$ws.Range("E41").Value = "'-0.23%"
$ws.Range("D42").Value = "'0.007340"
$ws.Range("E42").Value = "'-7.56%"
$ws.Range("D43").Value = "'0.001978"
$ws.Range("E43").Value = "'-1.95%"
$ws.Range("D44").Value = "'0.008237"
$ws.Range("E44").Value = "'-7.28%"
$ws.Range("D45").Value = "'0.00007084"
$ws.Range("E45").Value = "'-0.38%"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("D47").Value = "'0.0005799"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("E48").Value = "'54.18%"
$ws.Range("D49").Value = "'0.003512"
$ws.Range("E49").Value = "'23.60%"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E51").Value = "'0.04%"
